# Apply updates to Sheet1 per commit "Some successful tests w/ Brose."
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header labels
$ws.Range("B1").Value = "KvaserID"
$ws.Range("D1").Value = "MicrochipID"

# New data values for rows 2-6 (columns B, C, D, E)
$data = @(
    @("0x300", 122, "0x400", 119),
    @("0x203", 122, "0x401", 111),
    @("0x202", 12,  "0x402", 4),
    @("0x666", 5,   "0x403", 3),
    @("0x200", 2,   $null,   $null)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    if ($vals[2] -eq $null) {
        # Last data row (6) has no BikeID/Count entries; keep the cells
        # present but blank (matching the original empty-string cells)
        # without disturbing their default style.
        $ws.Cells.Item($row, 4).Style = "Normal"
    } else {
        $ws.Cells.Item($row, 4).Value = $vals[2]
    }
    if ($vals[3] -eq $null) {
        $ws.Cells.Item($row, 5).Style = "Normal"
    } else {
        $ws.Cells.Item($row, 5).Value = $vals[3]
    }
}

# Remove the now-obsolete rows 7-10 entirely (previously held items 5-8)
$ws.Range("A7:E10").Delete()
